$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: A7 text change ("High Cell" -> "Custom Height") + wrap text style
$ws.Range("A7").WrapText = $true
$ws.Range("A7").Value = "Custom Height"

# Row 1: A1 text change ("Wrap Text Row" -> "Initial Wrap Text ")
$ws.Range("A1").Value = "Initial Wrap Text "

# New row 9 data (order chosen to match shared-string insertion order)
$ws.Range("C9").Value = "looooong text"
$ws.Range("D9").Value = "veeeeeeeeeeeeeery loooooooooooong text"
$ws.Range("B9").Value = "loooong text"
$ws.Range("A9").Value = "Multiple"
$ws.Range("A9").Font.Bold = $true

# New column width for column A
$ws.Columns.Item(1).ColumnWidth = 9.5

$ws.Range("A10").Select() | Out-Null
